# Update the controlled-vocabulary wording for the two "reproductive structure"
# status descriptions (capitalize the first letter of each), and update the
# active-cell selection on the sheet from A2 to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$notYetObserved = "Reproductive/regenerative structures not yet observed since last fire"
$uncertainWhether = "Reproductive structures observed, but uncertain whether they were also produced earlier during post-fire regeneration"

# repr4 column (H:J) on the repr3a "uncertain/doubtful start" row
$ws.Range("H9").Value = $uncertainWhether
$ws.Range("I9").Value = $uncertainWhether
$ws.Range("J9").Value = $uncertainWhether

# repr3 / repr3a / repr4 columns (F:J) on the "no evidence" row
$ws.Range("F15").Value = $notYetObserved
$ws.Range("G15").Value = $notYetObserved
$ws.Range("H15").Value = $notYetObserved
$ws.Range("I15").Value = $notYetObserved
$ws.Range("J15").Value = $notYetObserved

# Move the active-cell selection from A2 to A3
$ws.Range("A3").Select() | Out-Null
